$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$players = @(
    "Keyonte George",
    "Lauri Markkanen",
    "Zach LaVine",
    "John Collins",
    "Jordan Poole",
    "Kyrie Irving",
    "CJ McCollum",
    "RJ Barrett",
    "Jimmy Butler",
    "Jalen Williams",
    "Brandon Boston Jr.",
    "Christian Braun",
    "Shai Gilgeous-Alexander",
    "Dennis Schröder",
    "Tobias Harris",
    "Joel Embiid"
)

$positions = @(
    "PG,SG",
    "SF,PF",
    "SG,SF",
    "PF,C",
    "PG,SG",
    "PG,SG",
    "PG,SG",
    "SF,PF",
    "SF,PF",
    "SG,SF,PF,C",
    "SG,SF,PF",
    "SG,SF",
    "PG",
    "PG",
    "SF,PF",
    "C"
)

$teams = @(
    "Utah Jazz",
    "Utah Jazz",
    "Chicago Bulls",
    "Utah Jazz",
    "Washington Wizards",
    "Dallas Mavericks",
    "New Orleans Pelicans",
    "Toronto Raptors",
    "Miami Heat",
    "Oklahoma City Thunder",
    "New Orleans Pelicans",
    "Denver Nuggets",
    "Oklahoma City Thunder",
    "Brooklyn Nets",
    "Detroit Pistons",
    "Philadelphia 76ers"
)

for ($i = 0; $i -lt $players.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $players[$i]
    $ws.Cells.Item($row, 2).Value = $positions[$i]
    $ws.Cells.Item($row, 3).Value = $teams[$i]
}
